# Apply the "graph step" update to the task/predecessor/value table.
#
# Before:
#   <task_name>   <predecessor>  <value>
#   <op1>                        10
#   <op2>         <op1>          20
#   <op3>         <op2>          30
#
# After:
#   <task_name>   <predecessor>  <value>
#   <op1>                        10
#   <op12>                       20
#   <op2-13>      <op1, op12>    30
#   <END>         <END>          <END>

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell text -------------------------------------------------
$ws.Range("A3").Value = "<op12>"
$ws.Range("B3").ClearContents()

$ws.Range("A4").Value = "<op2-13>"
$ws.Range("B4").Value = "<op1, op12>"

# --- Append the new terminating row --------------------------------------------
$ws.Range("A5").Value = "<END>"
$ws.Range("B5").Value = "<END>"
$ws.Range("C5").Value = "<END>"

# --- Formatting -----------------------------------------------------------------
# Re-apply the (same) Normal style to the "plain" cells so the workbook keeps a
# second cell format (mirrors the extra cellXfs entry produced upstream).
$ws.Range("A1:C1").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("C3").Style = "Normal"
$ws.Range("A4").Style = "Normal"
$ws.Range("C4").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Style = "Normal"

# Center the final "<END>" marker in the value column.
$ws.Range("C5").HorizontalAlignment = -4108

# --- Column widths ---------------------------------------------------------------
# Narrow the task/value columns back to the default width and widen the
# predecessor column so the longer "<op1, op12>" label fits.
$ws.Columns("A").ColumnWidth = 8.37
$ws.Columns("B").ColumnWidth = 10.97

# --- View / selection -------------------------------------------------------------
[void]$ws.Range("A2").Select()

$w = $excel.ActiveWindow
$w.TabRatio = 0.991
